# Apply "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across the resume body,
# matching the committed diff exactly: each metric substring becomes its
# own run with Bold + Color (RGB 0x2C3E50), splitting the original single
# run into multiple runs, with surrounding plain-text runs left untouched.

$d = $word.ActiveDocument

# RGB(0x2C, 0x3E, 0x50) expressed the way VBA's RGB() macro / Word COM's
# Font.Color integer expects: r + g*256 + b*65536
$metricColor = 44 + (62 * 256) + (80 * 65536)

function Highlight-Metrics {
    param(
        [string]$containsText,
        [string[]]$metrics
    )

    foreach ($p in $d.Paragraphs) {
        $pText = $p.Range.Text
        if ($pText.IndexOf($containsText) -ge 0) {
            $paraEnd = $p.Range.End
            $searchStart = $p.Range.Start
            foreach ($m in $metrics) {
                $rng = $d.Range($searchStart, $paraEnd)
                $found = $rng.Find.Execute($m, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
                if ($found) {
                    $rng.Font.Bold = $true
                    $rng.Font.Color = $metricColor
                    $searchStart = $rng.End
                }
            }
            break
        }
    }
}

Highlight-Metrics "Discovered systematic race coding errors" @("23%", "64%")

Highlight-Metrics "Utilized advanced sampling methods" @("±4.2%", "±2.1%", "71%", "87%")

Highlight-Metrics "Trigonometric algorithm for boundary estimation" @("73.5%", "$4.7M")

Highlight-Metrics "Built real-time FEC analysis systems" @("$2")

Highlight-Metrics "Modernized legacy ETL processes" @("57%")

Highlight-Metrics "Algorithmic innovation: Pioneered trigonometric boundary estimation" @("73.5%")

Highlight-Metrics "savings enabled nonprofit access" @("$4.7M")

Highlight-Metrics "Platform impact: Built redistricting system" @("12,847")

Write-Output "done"
